$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.529.05"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.914.71"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.80%  "
$ws.Range("D5").Value = "'326.02"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "'0.4827"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "'0.4072"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").Value = "'0.08145"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").Value = "'23.43"
$ws.Range("E11").Value = "  +3.61%  "
$ws.Range("D12").Value = "1.910.08"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'6.001"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'7.147"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "'90.20"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.06797"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "'0.00001040"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "29.539.77"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'5.625"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").Value = "'2.184"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("D25").Value = "2.200.91"
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("D26").Value = "'155.39"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "'6.418"
$ws.Range("E27").Value = "  +7.80%  "
$ws.Range("D28").Value = "'20.07"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").Value = "'119.73"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").Value = "'1.034"
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("D32").Value = "'0.09561"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'5.510"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "'3.566"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "'0.06103"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'1.175"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "'0.5936"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.986"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'10.73"
$ws.Range("E41").Value = "  +5.29%  "
$ws.Range("D42").Value = "'0.1857"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "'2.460"
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("D44").Value = "'1.281"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'0.07682"
$ws.Range("E45").Value = "  -4.02%  "
$ws.Range("D46").Value = "'12.42"
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").Value = "'0.5589"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "'1.943"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").Value = "'115.63"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").Value = "'72.67"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("E51").Value = "  +1.69%  "
